# Apply weekly crime-data refresh (NYPD 81st Precinct CompStat report)
# Volume 30, Number 32 -- week of 8/7/2023 through 8/13/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/issue number + reporting week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Cells whose type changes (number <-> placeholder text) are handled first,
#     copying cell format from an unchanged sibling cell so the workbook keeps
#     its existing "text" vs "numeric" style for these positions. ---
$ws.Range("G15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "'***.*"
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("D27").Value = 1
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = -100
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("C28").Value = "'0"
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("C29").Value = "'0"
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# --- Remaining crime statistics table updates ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = 21.428571428571
$ws.Range("L15").Value = 142.857142857143
$ws.Range("M15").Value = 6.25
$ws.Range("N15").Value = -65.306122448979
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 8.333333333333
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = 17.894736842105
$ws.Range("L16").Value = 13.131313131313
$ws.Range("M16").Value = -41.968911917098
$ws.Range("N16").Value = -84.594222833562
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -55.555555555555
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = -11.538461538461
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 217
$ws.Range("K17").Value = -8.294930875576
$ws.Range("L17").Value = 5.851063829787
$ws.Range("M17").Value = 1.015228426395
$ws.Range("N17").Value = -63.752276867031
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 100
$ws.Range("J18").Value = 122
$ws.Range("K18").Value = -18.032786885245
$ws.Range("L18").Value = 4.166666666666
$ws.Range("M18").Value = -21.875
$ws.Range("N18").Value = -78.902953586497
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -43.589743589743
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = -16.4
$ws.Range("L19").Value = -4.128440366972
$ws.Range("M19").Value = 18.75
$ws.Range("N19").Value = -8.333333333333
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 9
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = -25.882352941176
$ws.Range("M20").Value = 16.666666666666
$ws.Range("N20").Value = -83.846153846153
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -11.111111111111
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = -19.587628865979
$ws.Range("I21").Value = 702
$ws.Range("J21").Value = 788
$ws.Range("K21").Value = -10.913705583756
$ws.Range("L21").Value = 3.387334315169
$ws.Range("M21").Value = -9.536082474226
$ws.Range("N21").Value = -71.134868421052
$ws.Range("M22").Value = -37.5
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -14.285714285714
$ws.Range("I23").Value = 55
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = 5.76923076923
$ws.Range("L23").Value = -5.172413793103
$ws.Range("M23").Value = 7.843137254901
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = -14.851485148514
$ws.Range("I24").Value = 523
$ws.Range("J24").Value = 489
$ws.Range("K24").Value = 6.952965235173
$ws.Range("L24").Value = 49.856733524355
$ws.Range("M24").Value = 12.231759656652
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -5.714285714285
$ws.Range("I25").Value = 324
$ws.Range("J25").Value = 255
$ws.Range("K25").Value = 27.058823529411
$ws.Range("L25").Value = 55.76923076923
$ws.Range("M25").Value = -34.146341463414
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 23
$ws.Range("K26").Value = 15
$ws.Range("L26").Value = 76.923076923076
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -30
$ws.Range("L27").Value = -58.823529411764
$ws.Range("M28").Value = -72.972972972973
$ws.Range("M29").Value = -75
